# Reorder the MSME summary table so that "Enterprises density (per 1000 people)"
# (with its value "55") becomes the first data row of the table, right after the
# Micro/SMEs/MSMEs header row, pushing the other rows down by one.
#
# Final order (rows 12-16), column A = label, column D = value:
#   12: Enterprises density (per 1000 people) | 55
#   13: Employment (% of total)               | 86.2
#   14: Enterprises (absolute #)               | 1279784
#   15: Employment (absolute #)                | 8337000
#   16: Enterprises (% of total)               | 97.6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (labels) - plain text, safe to assign directly.
$ws.Range("A12").Value = "Enterprises density (per 1000 people)"
$ws.Range("A13").Value = "Employment (% of total)"
$ws.Range("A14").Value = "Enterprises (absolute #)"
$ws.Range("A15").Value = "Employment (absolute #)"
$ws.Range("A16").Value = "Enterprises (% of total)"

# Column D (values) - these are numeric-looking strings that must stay stored
# as text (matching the original file, which kept them as shared-string text,
# not real numbers). Temporarily force text format so Excel doesn't coerce
# them into numbers, then restore the original style/formatting afterwards.
$valRange = $ws.Range("D12:D16")
$origStyle = $valRange.Style

$valRange.NumberFormat = "@"

$ws.Range("D12").Value = "55"
$ws.Range("D13").Value = "86.2"
$ws.Range("D14").Value = "1279784"
$ws.Range("D15").Value = "8337000"
$ws.Range("D16").Value = "97.6"

$valRange.Style = $origStyle
